$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Desktop Computer" (DKS) device type rows (rows 8-10),
# which shifts all subsequent rows up.
$ws.Range("A8:G10").EntireRow.Delete()

# Set the active selection cell as recorded in the workbook after editing.
$ws.Range("E10").Select()

# Apply the page setup recorded for the sheet (paper size / orientation).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
